$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force a literal-text cell value even when $value looks like a number,
    # matching the source file (t="inlineStr"/shared-string, no numeric coercion),
    # then strip the resulting quote-prefix style so the cell's style index is
    # left exactly as it was (unstyled / General), matching the target XML.
    if ($value -match "^[+-]?[0-9]*\.?[0-9]+$") {
        $range.Value = "'" + $value
        $range.ClearFormats()
    } else {
        $range.Value = $value
    }
}

Set-TextValue $ws.Range("D2") '26.985.71'
Set-TextValue $ws.Range("E2") '  +0.20%  '
Set-TextValue $ws.Range("D3") '1.562.59'
Set-TextValue $ws.Range("E3") '  +0.58%  '
Set-TextValue $ws.Range("D4") '1.00'
Set-TextValue $ws.Range("E4") '  +0.07%  '
Set-TextValue $ws.Range("D5") '207.51'
Set-TextValue $ws.Range("E5") '  +0.33%  '
Set-TextValue $ws.Range("E6") '  +0.19%  '
Set-TextValue $ws.Range("E7") '  +0.04%  '
Set-TextValue $ws.Range("E8") '  +0.99%  '
Set-TextValue $ws.Range("E9") '  +0.39%  '
Set-TextValue $ws.Range("E10") '  +2.67%  '
Set-TextValue $ws.Range("E11") '  -0.38%  '
Set-TextValue $ws.Range("D12") '1.785.77'
Set-TextValue $ws.Range("E12") '  +0.52%  '
Set-TextValue $ws.Range("D13") '1.565.28'
Set-TextValue $ws.Range("E13") '  +0.60%  '
Set-TextValue $ws.Range("E14") '  +0.40%  '
Set-TextValue $ws.Range("E15") '  +0.54%  '
Set-TextValue $ws.Range("E16") '  +0.46%  '
Set-TextValue $ws.Range("D17") '26.974.11'
Set-TextValue $ws.Range("E17") '  +0.18%  '
Set-TextValue $ws.Range("D18") '0.0₃0706'
Set-TextValue $ws.Range("E18") '  +1.88%  '
Set-TextValue $ws.Range("D19") '215.97'
Set-TextValue $ws.Range("E19") '  -0.83%  '
Set-TextValue $ws.Range("D20") '7.36'
Set-TextValue $ws.Range("E20") '  +0.99%  '
Set-TextValue $ws.Range("D22") '4.10'
Set-TextValue $ws.Range("E22") '  +1.60%  '
Set-TextValue $ws.Range("D23") '9.20'
Set-TextValue $ws.Range("E23") '  +0.05%  '
Set-TextValue $ws.Range("E24") '  -1.14%  '
Set-TextValue $ws.Range("D25") '153.15'
Set-TextValue $ws.Range("E25") '  -0.68%  '
Set-TextValue $ws.Range("E26") '  +0.53%  '
Set-TextValue $ws.Range("D27") '15.08'
Set-TextValue $ws.Range("E27") '  +1.10%  '
Set-TextValue $ws.Range("E28") '  +1.61%  '
Set-TextValue $ws.Range("D29") '1.00'
Set-TextValue $ws.Range("E29") '  +0.00%  '
Set-TextValue $ws.Range("D30") '0.0471'
Set-TextValue $ws.Range("E30") '  +0.68%  '
Set-TextValue $ws.Range("E31") '  +1.55%  '
Set-TextValue $ws.Range("E32") '  +0.35%  '
Set-TextValue $ws.Range("D33") '3.12'
Set-TextValue $ws.Range("E33") '  +1.68%  '
Set-TextValue $ws.Range("D34") '1.421.38'
Set-TextValue $ws.Range("E34") '  -1.02%  '
Set-TextValue $ws.Range("D35") '1.61'
Set-TextValue $ws.Range("E35") '  +3.14%  '
Set-TextValue $ws.Range("E36") '  +8.77%  '
Set-TextValue $ws.Range("E37") '  +2.15%  '
Set-TextValue $ws.Range("E38") '  +0.26%  '
Set-TextValue $ws.Range("D39") '0.534'
Set-TextValue $ws.Range("E39") '  +2.79%  '
Set-TextValue $ws.Range("E40") '  +2.29%  '
Set-TextValue $ws.Range("E41") '  -0.40%  '
Set-TextValue $ws.Range("E42") '  -0.02%  '
Set-TextValue $ws.Range("E43") '  +2.47%  '
Set-TextValue $ws.Range("E44") '  +2.00%  '
Set-TextValue $ws.Range("E45") '  +0.94%  '
Set-TextValue $ws.Range("E46") '  -1.00%  '
Set-TextValue $ws.Range("D47") '1.699.39'
Set-TextValue $ws.Range("E47") '  +0.47%  '
Set-TextValue $ws.Range("D48") '87.36'
Set-TextValue $ws.Range("E48") '  +1.05%  '
Set-TextValue $ws.Range("E49") '  -0.83%  '
Set-TextValue $ws.Range("D50") '0.0₇0998'
Set-TextValue $ws.Range("E50") '  +2.05%  '
Set-TextValue $ws.Range("D51") '0.0959'
Set-TextValue $ws.Range("E51") '  +0.33%  '
